$wb = $excel.ActiveWorkbook

# --- FortiBranch sheet: update the firewall policy rows ---
$ws = $wb.Worksheets.Item("FortiBranch")

# Row 2: any/any/Deny -> all/Internal DB/http,https,mysql/Allow (VPN stays Enabled)
$ws.Range("A2").Value2 = "all"
$ws.Range("B2").Value2 = "Internal DB"
$ws.Range("C2").Value2 = "http,https,mysql"
$ws.Range("D2").Value2 = "Allow"

# Row 3: All_Internet/DMZ -> jb-m/webbranch (Services/Action/VPN unchanged)
$ws.Range("A3").Value2 = "jb-m"
$ws.Range("B3").Value2 = "webbranch"

# Row 4: ws12c/.../telnet... -> jb-m1/.../SSH (Destination/Action/VPN unchanged)
$ws.Range("A4").Value2 = "jb-m1"
$ws.Range("C4").Value2 = "SSH"

# Row 5: DMZ/Internal DB/telnet... -> w10-c/webbranch/SSH (Action/VPN unchanged)
$ws.Range("A5").Value2 = "w10-c"
$ws.Range("B5").Value2 = "webbranch"
$ws.Range("C5").Value2 = "SSH"

# Row 6 (webmain/webbranch/http,https/Allow/Enabled) is removed entirely
$ws.Rows.Item(6).Delete()

# --- CPMAN sheet: it is no longer the active tab, but keeps its own selection ---
$wsCpman = $wb.Worksheets.Item("CPMAN")
$wsCpman.Range("D30").Select()

# The FortiBranch tab becomes the selected/active one, with B7 selected
$ws.Activate()
$ws.Range("B7").Select()
